$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: SallieMae Loan (existing) gains AutoPay/Due Day/Amount ---
$ws.Range("B3").Value = "no"
$ws.Range("C3").Value = 28
$ws.Range("D3").Value = 103.24

# --- Row 4: ACS Loan (existing) gains Due Day/Amount, keeps AutoPay=Yes ---
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = 132.15

# --- Row 5: Sti Loan (existing) gains AutoPay/Due Day, keeps Amount ---
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = 26
$ws.Range("D5").Value = 222.16

# --- Row 7: Mortgage (new) ---
$ws.Range("A7").Value = "Mortgage"
$ws.Range("B7").Value = "no"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1588

# --- Row 8: water (new) ---
$ws.Range("A8").Value = "water"
$ws.Range("B8").Value = "no"
$ws.Range("D8").Value = 20

# --- Row 9: electric (new) ---
$ws.Range("A9").Value = "electric"
$ws.Range("B9").Value = "no"
$ws.Range("D9").Value = 150

# --- Row 11: Insurance (new) ---
$ws.Range("A11").Value = "Insurance"

# --- Row 12: mastercard (new) ---
$ws.Range("A12").Value = "mastercard"
$ws.Range("B12").Value = "no"

# --- Row 13: earthtreks (new) ---
$ws.Range("A13").Value = "earthtreks"
$ws.Range("B13").Value = "yes"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 70

# --- Row 20: total formula ---
$ws.Range("D20").Formula = "=SUM(D2:D19)"

# --- Selection matches the end-state cursor position ---
$ws.Range("D12").Select()
